# "add building (no collision detection)"
#
# The original deck had an ellipse shape labeled "Obstacle" in the legend
# on slide 5. This edit relabels it to "Building" (the obstacle/collision
# feature itself is not implemented yet - see commit message).

$p = $ppt.ActivePresentation

$targetSlideIndex = -1
$targetShapeIndex = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "Obstacle") {
                    $targetSlideIndex = $si
                    $targetShapeIndex = $i
                }
            }
        }
    }
}

if ($targetSlideIndex -eq -1) {
    throw "Could not locate the 'Obstacle' legend shape"
}

$slide = $p.Slides.Item($targetSlideIndex)
$shape = $slide.Shapes.Item($targetShapeIndex)

# Rename the legend entry; keeps existing run formatting (bold, size, etc.)
$shape.TextFrame.TextRange.Text = "Building"
